# Daily attendance processing - 2025-11-01 11:17:17
#
# The "Recorded By" column (G) lists the users who recorded / touched a
# session, separated by ", ". Whenever "System" was recorded as the FIRST
# entry in that list, move it to the LAST position instead (the rest of
# the list keeps its relative order). Cells where "System" is not the
# first entry (or is the only entry) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ', '

    if ($parts.Length -gt 1 -and $parts[0] -eq 'System') {
        $rest = $parts[1..($parts.Length - 1)]
        $newParts = $rest + ,$parts[0]
        $newText = [string]::Join(', ', $newParts)
        $cell.Value = $newText
    }
}
